$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 2428965.8
$ws.Range("I62").Value = 3482857.5
$ws.Range("J62").Value = 20070.5
$ws.Range("K62").Value = 3482857.5
$ws.Range("L62").Value = 20070.5
$ws.Range("M62").Value = -3482233.5
$ws.Range("N62").Value = -21318.5
# Row 64
$ws.Range("H64").Value = 4931.391
$ws.Range("I64").Value = 4030.8333
$ws.Range("J64").Value = 5913.8184
$ws.Range("K64").Value = 4030.8333
$ws.Range("L64").Value = 5913.8184
$ws.Range("M64").Value = -3782.8333
$ws.Range("N64").Value = -6409.8184
# Row 65
$ws.Range("H65").Value = 2428965.8
$ws.Range("I65").Value = 3482857.5
$ws.Range("J65").Value = 20070.5
$ws.Range("K65").Value = 17414287.5
$ws.Range("L65").Value = 100352.5
$ws.Range("M65").Value = -17411167.5
$ws.Range("N65").Value = -106592.5
# Row 67
$ws.Range("H67").Value = 4931.391
$ws.Range("I67").Value = 4030.8333
$ws.Range("J67").Value = 5913.8184
$ws.Range("K67").Value = 4030.8333
$ws.Range("L67").Value = 5913.8184
$ws.Range("M67").Value = -3172.8333
$ws.Range("N67").Value = -7629.8184
# Row 76
$ws.Range("H76").Value = 2648628
$ws.Range("I76").Value = 3475105
$ws.Range("J76").Value = 3902
$ws.Range("K76").Value = 3475105
$ws.Range("L76").Value = 3902
$ws.Range("M76").Value = -3474790
$ws.Range("N76").Value = -4532
# Row 79
$ws.Range("H79").Value = 2648628
$ws.Range("I79").Value = 3475105
$ws.Range("J79").Value = 3902
$ws.Range("K79").Value = 3475105
$ws.Range("L79").Value = 3902
$ws.Range("M79").Value = -3474013
$ws.Range("N79").Value = -6086
# Row 132
$ws.Range("H132").Value = 30484.543
$ws.Range("I132").Value = 33001.844
$ws.Range("J132").Value = 3633.3333
$ws.Range("K132").Value = 99005.53199999999
$ws.Range("L132").Value = 10899.9999
$ws.Range("M132").Value = -96475.53199999999
$ws.Range("N132").Value = -15959.9999
# Row 133
$ws.Range("H133").Value = 45267.895
$ws.Range("J133").Value = 45267.895
$ws.Range("L133").Value = 45267.895
$ws.Range("N133").Value = -55387.895
# Row 135
$ws.Range("H135").Value = 1690.1177
$ws.Range("I135").Value = 2120.9092
$ws.Range("J135").Value = 900.3333
$ws.Range("K135").Value = 19088.1828
$ws.Range("L135").Value = 8102.9997
$ws.Range("M135").Value = -16553.1828
$ws.Range("N135").Value = -13172.9997
# Row 138
$ws.Range("H138").Value = 6698913
$ws.Range("I138").Value = 1896215.5
$ws.Range("J138").Value = 9806541
$ws.Range("K138").Value = 5688646.5
$ws.Range("L138").Value = 29419623
$ws.Range("M138").Value = -5683506.5
$ws.Range("N138").Value = -29429903

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 25404.05
$ws.Range("I32").Value = 4988.659
$ws.Range("J32").Value = 81546.375
$ws.Range("K32").Value = 4988.659
$ws.Range("L32").Value = 81546.375
$ws.Range("M32").Value = -4701.659
$ws.Range("N32").Value = -82120.375
# Row 63
$ws.Range("H63").Value = 6011.88
$ws.Range("I63").Value = 4035.5
$ws.Range("J63").Value = 8527.272000000001
$ws.Range("K63").Value = 4035.5
$ws.Range("L63").Value = 8527.272000000001
$ws.Range("M63").Value = -3349.5
$ws.Range("N63").Value = -9899.272000000001
# Row 66
$ws.Range("H66").Value = 6011.88
$ws.Range("I66").Value = 4035.5
$ws.Range("J66").Value = 8527.272000000001
$ws.Range("K66").Value = 20177.5
$ws.Range("L66").Value = 42636.36
$ws.Range("M66").Value = -16745.5
$ws.Range("N66").Value = -49500.36
# Row 74
$ws.Range("H74").Value = 3608.7144
$ws.Range("I74").Value = 970.96875
$ws.Range("J74").Value = 8573.883
$ws.Range("K74").Value = 970.96875
$ws.Range("L74").Value = 8573.883
$ws.Range("M74").Value = -96.96875
$ws.Range("N74").Value = -10321.883
# Row 77
$ws.Range("H77").Value = 3608.7144
$ws.Range("I77").Value = 970.96875
$ws.Range("J77").Value = 8573.883
$ws.Range("K77").Value = 4854.84375
$ws.Range("L77").Value = 42869.415
$ws.Range("M77").Value = -486.84375
$ws.Range("N77").Value = -51605.415
# Row 88
$ws.Range("H88").Value = 1926.5
$ws.Range("I88").Value = 1884
$ws.Range("K88").Value = 1884
$ws.Range("M88").Value = -1478
# Row 91
$ws.Range("H91").Value = 1926.5
$ws.Range("I91").Value = 1884
$ws.Range("K91").Value = 1884
$ws.Range("M91").Value = -480
# Row 113
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = $null
# Row 133
$ws.Range("H133").Value = 45710.168
$ws.Range("J133").Value = 45710.168
$ws.Range("L133").Value = 45710.168
$ws.Range("N133").Value = -50770.168
# Row 139
$ws.Range("H139").Value = 45943
$ws.Range("J139").Value = 45943
$ws.Range("L139").Value = 45943
$ws.Range("N139").Value = -56223

$ws = $wb.Worksheets.Item("BSM")
# Row 133
$ws.Range("H133").Value = 43000
$ws.Range("J133").Value = 43000
$ws.Range("L133").Value = 43000
$ws.Range("N133").Value = -53120
# Row 134
$ws.Range("H134").Value = 3142.9697
$ws.Range("I134").Value = 2239.48
$ws.Range("J134").Value = 5966.375
$ws.Range("K134").Value = 6718.440000000001
$ws.Range("L134").Value = 17899.125
$ws.Range("M134").Value = -4183.440000000001
$ws.Range("N134").Value = -22969.125
# Row 139
$ws.Range("H139").Value = 85390
$ws.Range("J139").Value = 85390
$ws.Range("L139").Value = 85390
$ws.Range("N139").Value = -95670

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4533.6377
$ws.Range("I31").Value = 1301.3334
$ws.Range("J31").Value = 11168.368
$ws.Range("K31").Value = 1301.3334
$ws.Range("L31").Value = 11168.368
$ws.Range("M31").Value = -1006.3334
$ws.Range("N31").Value = -11758.368
# Row 34
$ws.Range("H34").Value = 4533.6377
$ws.Range("I34").Value = 1301.3334
$ws.Range("J34").Value = 11168.368
$ws.Range("K34").Value = 1301.3334
$ws.Range("L34").Value = 11168.368
$ws.Range("M34").Value = -1099.3334
$ws.Range("N34").Value = -11572.368
# Row 62
$ws.Range("H62").Value = 18639.643
$ws.Range("I62").Value = 21112.416
$ws.Range("J62").Value = 3803
$ws.Range("K62").Value = 21112.416
$ws.Range("L62").Value = 3803
$ws.Range("M62").Value = -20488.416
$ws.Range("N62").Value = -5051
# Row 65
$ws.Range("H65").Value = 18639.643
$ws.Range("I65").Value = 21112.416
$ws.Range("J65").Value = 3803
$ws.Range("K65").Value = 105562.08
$ws.Range("L65").Value = 19015
$ws.Range("M65").Value = -102442.08
$ws.Range("N65").Value = -25255
# Row 88
$ws.Range("H88").Value = 35000
$ws.Range("J88").Value = 35000
$ws.Range("L88").Value = 35000
$ws.Range("N88").Value = -35812
# Row 91
$ws.Range("H91").Value = 35000
$ws.Range("J91").Value = 35000
$ws.Range("L91").Value = 35000
$ws.Range("N91").Value = -37808

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1894.1
$ws.Range("I5").Value = 1258.5294
$ws.Range("J5").Value = 2221.5151
$ws.Range("K5").Value = 3775.5882
$ws.Range("L5").Value = 6664.5453
$ws.Range("M5").Value = -3663.5882
$ws.Range("N5").Value = -6888.5453
# Row 107
$ws.Range("H107").Value = 545.8
$ws.Range("I107").Value = 582.6667
$ws.Range("J107").Value = 490.5
$ws.Range("K107").Value = 1748.0001
$ws.Range("L107").Value = 1471.5
$ws.Range("M107").Value = 171.9999
$ws.Range("N107").Value = -5311.5
# Row 135
$ws.Range("H135").Value = 1894.1
$ws.Range("I135").Value = 1258.5294
$ws.Range("J135").Value = 2221.5151
$ws.Range("K135").Value = 11326.7646
$ws.Range("L135").Value = 19993.6359
$ws.Range("M135").Value = -8791.764599999999
$ws.Range("N135").Value = -25063.6359

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3179.2856
$ws.Range("I80").Value = 3029.0625
$ws.Range("J80").Value = 3660
$ws.Range("K80").Value = 3029.0625
$ws.Range("L80").Value = 3660
$ws.Range("M80").Value = -2031.0625
$ws.Range("N80").Value = -5656
# Row 83
$ws.Range("H83").Value = 3179.2856
$ws.Range("I83").Value = 3029.0625
$ws.Range("J83").Value = 3660
$ws.Range("K83").Value = 15145.3125
$ws.Range("L83").Value = 18300
$ws.Range("M83").Value = -10153.3125
$ws.Range("N83").Value = -28284
# Row 113
$ws.Range("H113").Value = 1181.6666
$ws.Range("I113").Value = 1038
$ws.Range("J113").Value = 1296.6
$ws.Range("K113").Value = 1038
$ws.Range("L113").Value = 1296.6
$ws.Range("M113").Value = 1132
$ws.Range("N113").Value = -5636.6
# Row 126
$ws.Range("H126").Value = 2519.261
$ws.Range("I126").Value = 2089.8096
$ws.Range("K126").Value = 6269.4288
$ws.Range("M126").Value = -3799.4288
# Row 138
$ws.Range("H138").Value = 75485.8
$ws.Range("J138").Value = 75485.8
$ws.Range("L138").Value = 75485.8
$ws.Range("N138").Value = -85765.8
# Row 139
$ws.Range("H139").Value = 43403.777
$ws.Range("J139").Value = 43403.777
$ws.Range("L139").Value = 43403.777
$ws.Range("N139").Value = -53683.777

$ws = $wb.Worksheets.Item("WVR")
# Row 51
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").Value = $null
# Row 107
$ws.Range("H107").Value = 618.46155
$ws.Range("I107").Value = 970
$ws.Range("J107").Value = 398.75
$ws.Range("K107").Value = 2910
$ws.Range("L107").Value = 1196.25
$ws.Range("M107").Value = -990
$ws.Range("N107").Value = -5036.25
# Row 136
$ws.Range("H136").Value = 1992.8
$ws.Range("I136").Value = 1234.2
$ws.Range("J136").Value = 3510
$ws.Range("K136").Value = 3702.6
$ws.Range("L136").Value = 10530
$ws.Range("M136").Value = -1152.6
$ws.Range("N136").Value = -15630
